$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New narrower column for the "Name" label column (K)
$ws.Columns.Item(11).ColumnWidth = 9.3

# New "Name" header above the new grain-size/flux/diffusivity table
$ws.Range("J25").Value = "Name"

# New header row for the grain-size / flux / diffusivity table
$ws.Range("I26").Value = "Grain size"
$ws.Range("J26").Value = "Flux"
$ws.Range("K26").Value = "Diffusivity"

# Remove the old footnote that used to live next to the F27 total
$ws.Range("G27").ClearContents()

# New diffusivity row under row 27 (divide each flux figure by the 0.244 g/cm2 factor)
$ws.Range("A28").Formula = "=A27/0.244"
$ws.Range("B28").Formula = "=B27/0.244"
$ws.Range("C28").Formula = "=C27/0.244"
$ws.Range("D28").Formula = "=D27/0.244"
$ws.Range("E28").Formula = "=E27/0.244"
$ws.Range("A28:E28").Style = $ws.Range("A27").Style

# New diffusivity row under row 31 (divide each flux figure by the 0.349 g/cm2 factor)
$ws.Range("A32").Formula = "=A31/0.349"
$ws.Range("B32").Formula = "=B31/0.349"
$ws.Range("C32").Formula = "=C31/0.349"
$ws.Range("D32").Formula = "=D31/0.349"
$ws.Range("E32").Formula = "=E31/0.349"
$ws.Range("A32:E32").Style = $ws.Range("A31").Style

# Update selection to match the last-edited cell
$ws.Range("E32").Select()
